# Update evaluation metrics across the three result sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.3950177935943061
$ws1.Range("C2").Value = 0.07608695652173914
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 0.1414141414141414
$ws1.Range("F2").Value = 0.2916666666666667
$ws1.Range("G2").Value = 0.6816479400749064
$ws1.Range("H2").Value = 0.7891252006420545
$ws1.Range("I2").Value = 28
$ws1.Range("J2").Value = 340
$ws1.Range("K2").Value = 194
$ws1.Range("L2").Value = 0

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 0.3632958801498127
$ws2.Range("D2").Value = 0.532967032967033

$ws2.Range("B3").Value = 0.07608695652173914
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.1414141414141414

$ws2.Range("B4").Value = 0.3950177935943061
$ws2.Range("C4").Value = 0.3950177935943061
$ws2.Range("D4").Value = 0.3950177935943061
$ws2.Range("E4").Value = 0.3950177935943061

$ws2.Range("B5").Value = 0.5380434782608696
$ws2.Range("C5").Value = 0.6816479400749064
$ws2.Range("D5").Value = 0.3371905871905872

$ws2.Range("B6").Value = 0.9539687451647842
$ws2.Range("C6").Value = 0.3950177935943061
$ws2.Range("D6").Value = 0.5134590597224049

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 194
$ws3.Range("C2").Value = 340
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 28
